$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "62.423.94"
$cell = $ws.Range("E2")
$cell.NumberFormat = "@"
$cell.Value = "  -2.24%  "
$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "3.037.13"
$cell = $ws.Range("E3")
$cell.NumberFormat = "@"
$cell.Value = "  -2.54%  "
$cell = $ws.Range("E4")
$cell.NumberFormat = "@"
$cell.Value = "  +0.28%  "
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "532.46"
$cell = $ws.Range("E5")
$cell.NumberFormat = "@"
$cell.Value = "  -4.82%  "
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "131.80"
$cell = $ws.Range("E6")
$cell.NumberFormat = "@"
$cell.Value = "  -5.58%  "
$cell = $ws.Range("E7")
$cell.NumberFormat = "@"
$cell.Value = "  +0.08%  "
$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "3.032.31"
$cell = $ws.Range("E8")
$cell.NumberFormat = "@"
$cell.Value = "  -2.49%  "
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "0.494"
$cell = $ws.Range("E9")
$cell.NumberFormat = "@"
$cell.Value = "  +0.09%  "
$cell = $ws.Range("E10")
$cell.NumberFormat = "@"
$cell.Value = "  -1.28%  "
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "6.14"
$cell = $ws.Range("E11")
$cell.NumberFormat = "@"
$cell.Value = "  -9.46%  "
$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "0.450"
$cell = $ws.Range("E12")
$cell.NumberFormat = "@"
$cell.Value = "  -1.48%  "
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "0.0000222"
$cell = $ws.Range("E13")
$cell.NumberFormat = "@"
$cell.Value = "  +1.92%  "
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "33.92"
$cell = $ws.Range("E14")
$cell.NumberFormat = "@"
$cell.Value = "  -4.96%  "
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "3.536.24"
$cell = $ws.Range("E15")
$cell.NumberFormat = "@"
$cell.Value = "  -2.18%  "
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "62.470.07"
$cell = $ws.Range("E16")
$cell.NumberFormat = "@"
$cell.Value = "  -2.02%  "
$cell = $ws.Range("E17")
$cell.NumberFormat = "@"
$cell.Value = "  -1.11%  "
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "3.046.17"
$cell = $ws.Range("E18")
$cell.NumberFormat = "@"
$cell.Value = "  -1.88%  "
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "6.56"
$cell = $ws.Range("E19")
$cell.NumberFormat = "@"
$cell.Value = "  -2.02%  "
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "479.17"
$cell = $ws.Range("E20")
$cell.NumberFormat = "@"
$cell.Value = "  -5.51%  "
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "13.15"
$cell = $ws.Range("E21")
$cell.NumberFormat = "@"
$cell.Value = "  -4.67%  "
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "0.688"
$cell = $ws.Range("E22")
$cell.NumberFormat = "@"
$cell.Value = "  -2.93%  "
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "7.05"
$cell = $ws.Range("E23")
$cell.NumberFormat = "@"
$cell.Value = "  -3.15%  "
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "78.76"
$cell = $ws.Range("E24")
$cell.NumberFormat = "@"
$cell.Value = "  +1.04%  "
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "12.00"
$cell = $ws.Range("E25")
$cell.NumberFormat = "@"
$cell.Value = "  -3.75%  "
$cell = $ws.Range("E26")
$cell.NumberFormat = "@"
$cell.Value = "  -0.04%  "
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "2.67"
$cell = $ws.Range("E27")
$cell.NumberFormat = "@"
$cell.Value = "  -4.05%  "
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "7.99"
$cell = $ws.Range("E28")
$cell.NumberFormat = "@"
$cell.Value = "  -5.44%  "
$cell = $ws.Range("E29")
$cell.NumberFormat = "@"
$cell.Value = "  +0.46%  "
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "25.71"
$cell = $ws.Range("E30")
$cell.NumberFormat = "@"
$cell.Value = "  -2.53%  "
$cell = $ws.Range("E31")
$cell.NumberFormat = "@"
$cell.Value = "  -10.14%  "
$cell = $ws.Range("E32")
$cell.NumberFormat = "@"
$cell.Value = "  -1.26%  "
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "2.34"
$cell = $ws.Range("E33")
$cell.NumberFormat = "@"
$cell.Value = "  -8.80%  "
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "56.38"
$cell = $ws.Range("E34")
$cell.NumberFormat = "@"
$cell.Value = "  +0.61%  "
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "5.32"
$cell = $ws.Range("E35")
$cell.NumberFormat = "@"
$cell.Value = "  +1.80%  "
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "5.89"
$cell = $ws.Range("E36")
$cell.NumberFormat = "@"
$cell.Value = "  -0.69%  "
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "472.88"
$cell = $ws.Range("E37")
$cell.NumberFormat = "@"
$cell.Value = "  -13.20%  "
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "0.0392"
$cell = $ws.Range("E38")
$cell.NumberFormat = "@"
$cell.Value = "  -6.03%  "
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "3.069.01"
$cell = $ws.Range("E39")
$cell.NumberFormat = "@"
$cell.Value = "  -0.35%  "
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "0.0790"
$cell = $ws.Range("E40")
$cell.NumberFormat = "@"
$cell.Value = "  -1.61%  "
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "0.114"
$cell = $ws.Range("E41")
$cell.NumberFormat = "@"
$cell.Value = "  -4.41%  "
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "8.03"
$cell = $ws.Range("E42")
$cell.NumberFormat = "@"
$cell.Value = "  -1.47%  "
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "2.62"
$cell = $ws.Range("E43")
$cell.NumberFormat = "@"
$cell.Value = "  -0.14%  "
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "0.250"
$cell = $ws.Range("E44")
$cell.NumberFormat = "@"
$cell.Value = "  -2.43%  "
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "0.0₃0541"
$cell = $ws.Range("E46")
$cell.NumberFormat = "@"
$cell.Value = "  +7.20%  "
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "2.01"
$cell = $ws.Range("E47")
$cell.NumberFormat = "@"
$cell.Value = "  -4.80%  "
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "119.93"
$cell = $ws.Range("E48")
$cell.NumberFormat = "@"
$cell.Value = "  -1.82%  "
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "24.34"
$cell = $ws.Range("E49")
$cell.NumberFormat = "@"
$cell.Value = "  -0.66%  "
$cell = $ws.Range("E50")
$cell.NumberFormat = "@"
$cell.Value = "  +0.14%  "
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "2.31"
$cell = $ws.Range("E51")
$cell.NumberFormat = "@"
$cell.Value = "  +1.99%  "
